$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.5551036666666667
$ws.Range("H2").Value = 1.665311
$ws.Range("I2").Value = 0.389144998960137
$ws.Range("J2").Value = 0.389144998960137
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.7328106666666666
$ws.Range("N2").Value = 2.198432
$ws.Range("O2").Value = 0.4072614640191846
$ws.Range("P2").Value = 0.4072614640191846
$ws.Range("Q2").Value = 0.4067858880391111
$ws.Range("R2").Value = 3.661072992352
$ws.Range("S2").Value = 0.1584837619922495
$ws.Range("T2").Value = 0.1584837619922495

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.5551036666666667
$ws.Range("H3").Value = 1.665311
$ws.Range("I3").Value = 0.389144998960137
$ws.Range("J3").Value = 0.389144998960137
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.6260680000000001
$ws.Range("N3").Value = 1.878204
$ws.Range("O3").Value = 0.3479389450147599
$ws.Range("P3").Value = 0.3479389450147599
$ws.Range("Q3").Value = 0.3475326423826667
$ws.Range("R3").Value = 3.127793781444
$ws.Range("S3").Value = 0.1353987003959599
$ws.Range("T3").Value = 0.1353987003959599

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.5551036666666667
$ws.Range("H4").Value = 1.665311
$ws.Range("I4").Value = 0.389144998960137
$ws.Range("J4").Value = 0.389144998960137
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.440483
$ws.Range("N4").Value = 1.321449
$ws.Range("O4").Value = 0.2447995909660556
$ws.Range("P4").Value = 0.2447995909660555
$ws.Range("Q4").Value = 0.2445137284043333
$ws.Range("R4").Value = 2.200623555639
$ws.Range("S4").Value = 0.09526253657192767
$ws.Range("T4").Value = 0.09526253657192764

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.6252976666666666
$ws.Range("H5").Value = 1.875893
$ws.Range("I5").Value = 0.4383531842006258
$ws.Range("J5").Value = 0.4383531842006257
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.7328106666666666
$ws.Range("N5").Value = 2.198432
$ws.Range("O5").Value = 0.4072614640191846
$ws.Range("P5").Value = 0.4072614640191846
$ws.Range("Q5").Value = 0.4582247999751111
$ws.Range("R5").Value = 4.124023199776
$ws.Range("S5").Value = 0.1785243595550181
$ws.Range("T5").Value = 0.1785243595550181

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.6252976666666666
$ws.Range("H6").Value = 1.875893
$ws.Range("I6").Value = 0.4383531842006258
$ws.Range("J6").Value = 0.4383531842006257
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.6260680000000001
$ws.Range("N6").Value = 1.878204
$ws.Range("O6").Value = 0.3479389450147599
$ws.Range("P6").Value = 0.3479389450147599
$ws.Range("Q6").Value = 0.3914788595746667
$ws.Range("R6").Value = 3.523309736172
$ws.Range("S6").Value = 0.1525201444546265
$ws.Range("T6").Value = 0.1525201444546264

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.6252976666666666
$ws.Range("H7").Value = 1.875893
$ws.Range("I7").Value = 0.4383531842006258
$ws.Range("J7").Value = 0.4383531842006257
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.440483
$ws.Range("N7").Value = 1.321449
$ws.Range("O7").Value = 0.2447995909660556
$ws.Range("P7").Value = 0.2447995909660555
$ws.Range("Q7").Value = 0.2754329921063333
$ws.Range("R7").Value = 2.478896928957
$ws.Range("S7").Value = 0.1073086801909812
$ws.Range("T7").Value = 0.1073086801909812

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.2460686666666667
$ws.Range("H8").Value = 0.738206
$ws.Range("I8").Value = 0.1725018168392372
$ws.Range("J8").Value = 0.1725018168392372
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.7328106666666666
$ws.Range("N8").Value = 2.198432
$ws.Range("O8").Value = 0.4072614640191846
$ws.Range("P8").Value = 0.4072614640191846
$ws.Range("Q8").Value = 0.1803217436657778
$ws.Range("R8").Value = 1.622895692992
$ws.Range("S8").Value = 0.07025334247191697
$ws.Range("T8").Value = 0.07025334247191696

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.2460686666666667
$ws.Range("H9").Value = 0.738206
$ws.Range("I9").Value = 0.1725018168392372
$ws.Range("J9").Value = 0.1725018168392372
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.6260680000000001
$ws.Range("N9").Value = 1.878204
$ws.Range("O9").Value = 0.3479389450147599
$ws.Range("P9").Value = 0.3479389450147599
$ws.Range("Q9").Value = 0.1540557180026667
$ws.Range("R9").Value = 1.386501462024
$ws.Range("S9").Value = 0.06002010016417354
$ws.Range("T9").Value = 0.06002010016417353

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.2460686666666667
$ws.Range("H10").Value = 0.738206
$ws.Range("I10").Value = 0.1725018168392372
$ws.Range("J10").Value = 0.1725018168392372
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.440483
$ws.Range("N10").Value = 1.321449
$ws.Range("O10").Value = 0.2447995909660556
$ws.Range("P10").Value = 0.2447995909660555
$ws.Range("Q10").Value = 0.1083890644993333
$ws.Range("R10").Value = 0.9755015804939999
$ws.Range("S10").Value = 0.0422283742031467
$ws.Range("T10").Value = 0.04222837420314669

